# S23/G05: Trade dialog execution routing (manual/auto + live/paper)
# Updates row heights / remarks-column alignment for rows 186-194 and
# appends four new task rows (195-198) describing the new G05 work item.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Tweak row heights for the existing S23/G01-G04 rows (186-194) and
#    normalise the "remarks" (column H) cell alignment to match the
#    wrap-top style used by the rest of the table.
# ---------------------------------------------------------------------
$newHeights = @{
    186 = 27.75
    187 = 41.25
    188 = 27.75
    189 = 27.75
    190 = 27.75
    191 = 27.75
    192 = 27.75
    193 = 27.75
    194 = 27.75
}

foreach ($r in 186..194) {
    $ws.Rows.Item($r).RowHeight = $newHeights[$r]

    $hCell = $ws.Range("H$r")
    $hCell.WrapText = $true
    $hCell.VerticalAlignment = -4160   # xlVAlignTop
}

# ---------------------------------------------------------------------
# 2. Append the new S23/G05 rows (195-198).
# ---------------------------------------------------------------------
$newRows = @(
    @{
        Row = 195
        Height = 41.75
        A = "S23"
        B = "G05"
        C = "Trade dialog execution routing (manual/auto + live/paper)"
        D = "S23_G05_TB001"
        E = "Add per-order execution_target (LIVE/PAPER) and extend order create/execute flow to support AUTO (send now) vs MANUAL (queue) with audit trail."
        G = "implemented"
        H = "Supports trade dialogs; auto skips Waiting Queue but persists in Orders history."
    },
    @{
        Row = 196
        Height = 41.75
        A = "S23"
        B = "G05"
        C = "Trade dialog execution routing (manual/auto + live/paper)"
        D = "S23_G05_TF001"
        E = "Add Mode (MANUAL/AUTO) + Execution Target (LIVE/PAPER) controls to per-symbol and bulk buy/sell dialogs; default MANUAL + LIVE; confirm on AUTO+LIVE."
        G = "implemented"
        H = "AUTO executes immediately; MANUAL enqueues. Bulk executes sequentially with progress."
    },
    @{
        Row = 197
        Height = 28.35
        A = "S23"
        B = "G05"
        C = "Trade dialog execution routing (manual/auto + live/paper)"
        D = "S23_G05_TF002"
        E = "Add columns to show order mode + execution_target in Orders/Queue UIs and enable quick filtering."
        G = "implemented"
        H = "Queue primarily shows MANUAL, but mode column helps audit and debugging."
    },
    @{
        Row = 198
        Height = 41.75
        A = "S23"
        B = "G05"
        C = "Trade dialog execution routing (manual/auto + live/paper)"
        D = "S23_G05_TF003"
        E = "Add Settings defaults for trade dialogs (default mode + default execution target) and persist in localStorage/server config as appropriate."
        G = "planned"
        H = "Keeps UX safe by default; user can opt-in to AUTO/LIVE."
    }
)

foreach ($row in $newRows) {
    $r = $row.Row

    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H

    foreach ($col in @("A", "B", "C", "D", "E", "G", "H")) {
        $cell = $ws.Range("$col$r")
        $cell.WrapText = $true
        $cell.VerticalAlignment = -4160   # xlVAlignTop
    }

    $ws.Rows.Item($r).RowHeight = $row.Height
}

# ---------------------------------------------------------------------
# 3. Restore the view/selection state (scroll position + active cell).
# ---------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 185
$win.ScrollColumn = 1
$ws.Range("B195").Select()
